$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.035.24'
$ws.Range("E2").Value = '  -0.36%  '

$ws.Range("D3").Value = '''1.872.31'
$ws.Range("E3").Value = '  -2.61%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '''319.41'

$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").Value = '''0.5038'
$ws.Range("E7").Value = '  -3.50%  '

$ws.Range("D8").Value = '''0.3956'
$ws.Range("E8").Value = '  -3.49%  '

$ws.Range("D9").Value = '''0.08196'
$ws.Range("E9").Value = '  -4.05%  '

$ws.Range("D10").Value = '''42.21'
$ws.Range("E10").Value = '  -2.85%  '

$ws.Range("D11").Value = '''1.092'
$ws.Range("E11").Value = '  -3.38%  '

$ws.Range("D12").Value = '''23.71'
$ws.Range("E12").Value = '  +5.42%  '

$ws.Range("D13").Value = '''1.862.17'
$ws.Range("E13").Value = '  -3.21%  '

$ws.Range("D14").Value = '''6.291'
$ws.Range("E14").Value = '  -2.12%  '

$ws.Range("D15").Value = '''7.185'
$ws.Range("E15").Value = '  -3.45%  '

$ws.Range("E16").Value = '  +0.14%  '

$ws.Range("D17").Value = '''91.82'
$ws.Range("E17").Value = '  -4.50%  '

$ws.Range("D18").Value = '''0.00001089'
$ws.Range("E18").Value = '  -2.42%  '

$ws.Range("D19").Value = '''0.06403'
$ws.Range("E19").Value = '  -4.68%  '

$ws.Range("D20").Value = '''18.11'
$ws.Range("E20").Value = '  -1.28%  '

$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").Value = '''30.036.59'
$ws.Range("E22").Value = '  -0.39%  '

$ws.Range("D23").Value = '''5.839'
$ws.Range("E23").Value = '  -3.46%  '

$ws.Range("D24").Value = '''11.13'
$ws.Range("E24").Value = '  -1.90%  '

$ws.Range("D25").Value = '''2.175'
$ws.Range("E25").Value = '  -2.09%  '

$ws.Range("D26").Value = '''2.083.22'
$ws.Range("E26").Value = '  -2.85%  '

$ws.Range("D27").Value = '''21.42'
$ws.Range("E27").Value = '  +1.10%  '

$ws.Range("D28").Value = '''160.26'
$ws.Range("E28").Value = '  +0.34%  '

$ws.Range("D29").Value = '''2.228'
$ws.Range("E29").Value = '  -9.61%  '

$ws.Range("D30").Value = '''127.19'
$ws.Range("E30").Value = '  -1.69%  '

$ws.Range("D31").Value = '''1.066'
$ws.Range("E31").Value = '  -1.66%  '

$ws.Range("D32").Value = '''0.1033'

$ws.Range("D33").Value = '''5.935'
$ws.Range("E33").Value = '  -2.99%  '

$ws.Range("D34").Value = '''3.685'
$ws.Range("E34").Value = '  +1.15%  '

$ws.Range("D35").Value = '''0.02439'
$ws.Range("E35").Value = '  -2.94%  '

$ws.Range("D36").Value = '''5.231'
$ws.Range("E36").Value = '  -0.26%  '

$ws.Range("D37").Value = '''0.06368'
$ws.Range("E37").Value = '  -3.68%  '

$ws.Range("E38").Value = '  -3.41%  '

$ws.Range("D39").Value = '''1.174'
$ws.Range("E39").Value = '  -5.23%  '

$ws.Range("D40").Value = '''8.499'
$ws.Range("E40").Value = '  -5.22%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.6307'
$ws.Range("E41").Value = '  -3.66%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''1.216'
$ws.Range("E42").Value = '  -2.43%  '

$ws.Range("D43").Value = '''11.31'
$ws.Range("E43").Value = '  -3.87%  '

$ws.Range("D44").Value = '''1.0000'
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''12.98'
$ws.Range("E45").Value = '  -2.13%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.5913'
$ws.Range("E46").Value = '  -4.41%  '

$ws.Range("D47").Value = '''2.107'
$ws.Range("E47").Value = '  +0.51%  '

$ws.Range("E48").Value = '  -3.93%  '

$ws.Range("D49").Value = '''122.99'
$ws.Range("E49").Value = '  -1.70%  '

$ws.Range("E50").Value = '  -3.49%  '

$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").Value = '''1.124'
$ws.Range("E51").Value = '  -3.47%  '
